# Update the "ultimas semanas" (last weeks) actuals on the TASK sheet.
# The "Código" task row (row 30) had its actual-hours figure (column K)
# revised from 20 to 5.4; the running-total formula in column L
# recalculates automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("K30").Value = 5.4
